# Generate Report for Handoff
# Rewrites the Overview / zh-cn / de-de sheets so that the four tracked
# e2e files are renamed from GUID-based file names to the friendlier
# calleeMd1 / calleeMd2 / callerMd1 / callerMd2 names, refreshes their
# handoff timestamps + handoff xlf file names, fixes a couple of
# status/extension/dependency cells, and appends a 4th data row
# (callerMd2) to every sheet.

$wb = $excel.ActiveWorkbook

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Values (A2:D4 renamed, row 5 added)
$ws.Range("A2").Value = "calleeMd1.md"
$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("D2").Value = "2016-03-21 19:01:30"

$ws.Range("A3").Value = "calleeMd2.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "2016-03-21 19:01:30"

$ws.Range("A4").Value = "callerMd1.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("D4").Value = "2016-03-21 19:01:30"

$ws.Range("A5").Value = "callerMd2.md"
$ws.Range("B5").Value = "Ready for handoff"
$ws.Range("C5").Value = "Ready for handoff"
$ws.Range("D5").Value = "2016-03-21 19:01:30"
$ws.Range("D5").NumberFormat = $dateFmt

# Hyperlinks: drop all existing ones and recreate them (the engine's
# partial-hyperlink delete is unreliable), preserving A2..A4 and adding A5.
$ws.Hyperlinks.Delete()
$base = "https://github.com/OpenLocalizationTest/oltest/blob/d21edc54a3494f898cf9a13269ff04c61eedb23a/e2e/"
$ws.Hyperlinks.Add($ws.Range("A2"), ($base + "calleeMd1.md"), [Type]::Missing, [Type]::Missing, "calleeMd1.md")
$ws.Hyperlinks.Add($ws.Range("A3"), ($base + "calleeMd2.md"), [Type]::Missing, [Type]::Missing, "calleeMd2.md")
$ws.Hyperlinks.Add($ws.Range("A4"), ($base + "callerMd1.md"), [Type]::Missing, [Type]::Missing, "callerMd1.md")
$ws.Hyperlinks.Add($ws.Range("A5"), ($base + "callerMd2.md"), [Type]::Missing, [Type]::Missing, "callerMd2.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "calleeMd1.md"
$ws.Range("D2").Value = "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.zh-cn.xlf"
$ws.Range("E2").Value = "2016-03-21 19:01:22"
$ws.Range("K2").Value = "e2e\callerMd2.md,`ne2e\callerMd1.md"

$ws.Range("A3").Value = "calleeMd2.md"
$ws.Range("B3").Value = ".md"
$ws.Range("D3").Value = "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.zh-cn.xlf"
$ws.Range("E3").Value = "2016-03-21 19:01:22"
$ws.Range("J3").Value = "Include"
$ws.Range("K3").Value = "e2e\callerMd1.md"

$ws.Range("A4").Value = "callerMd1.md"
$ws.Range("B4").Value = ".md"
$ws.Range("D4").Value = "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.zh-cn.xlf"
$ws.Range("E4").Value = "2016-03-21 19:01:22"
$ws.Range("I4").Value = "e2e\calleeMd1.md,`ne2e\calleeMd2.md"
$ws.Range("J4").Value = "Include"
$ws.Range("K4").ClearContents()

$ws.Range("A5").Value = "callerMd2.md"
$ws.Range("B5").Value = ".md"
$ws.Range("C5").Value = "Ready for handoff"
$ws.Range("D5").Value = "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.zh-cn.xlf"
$ws.Range("E5").Value = "2016-03-21 19:01:22"
$ws.Range("H5").Value = "0001-01-01 00:00:00"
$ws.Range("I5").Value = "e2e\calleeMd1.md"
$ws.Range("J5").Value = "Include"
$ws.Range("E5").NumberFormat = $dateFmt
$ws.Range("H5").NumberFormat = $dateFmt

$ws.Hyperlinks.Delete()
$baseMd = "https://github.com/OpenLocalizationTest/oltest/blob/d21edc54a3494f898cf9a13269ff04c61eedb23a/e2e/"
$baseXlf = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9d673cf605f35295340e002a4ba6a1bb93c92ee0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/"
$ws.Hyperlinks.Add($ws.Range("A2"), ($baseMd + "calleeMd1.md"), [Type]::Missing, [Type]::Missing, "calleeMd1.md")
$ws.Hyperlinks.Add($ws.Range("D2"), ($baseXlf + "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.zh-cn.xlf"), [Type]::Missing, [Type]::Missing, "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), ($baseMd + "calleeMd2.md"), [Type]::Missing, [Type]::Missing, "calleeMd2.md")
$ws.Hyperlinks.Add($ws.Range("D3"), ($baseXlf + "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.zh-cn.xlf"), [Type]::Missing, [Type]::Missing, "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A4"), ($baseMd + "callerMd1.md"), [Type]::Missing, [Type]::Missing, "callerMd1.md")
$ws.Hyperlinks.Add($ws.Range("D4"), ($baseXlf + "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.zh-cn.xlf"), [Type]::Missing, [Type]::Missing, "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A5"), ($baseMd + "callerMd2.md"), [Type]::Missing, [Type]::Missing, "callerMd2.md")
$ws.Hyperlinks.Add($ws.Range("D5"), ($baseXlf + "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.zh-cn.xlf"), [Type]::Missing, [Type]::Missing, "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.zh-cn.xlf")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "calleeMd1.md"
$ws.Range("D2").Value = "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.de-de.xlf"
$ws.Range("E2").Value = "2016-03-21 19:01:30"
$ws.Range("K2").Value = "e2e\callerMd2.md,`ne2e\callerMd1.md"

$ws.Range("A3").Value = "calleeMd2.md"
$ws.Range("B3").Value = ".md"
$ws.Range("D3").Value = "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.de-de.xlf"
$ws.Range("E3").Value = "2016-03-21 19:01:30"
$ws.Range("J3").Value = "Include"
$ws.Range("K3").Value = "e2e\callerMd1.md"

$ws.Range("A4").Value = "callerMd1.md"
$ws.Range("B4").Value = ".md"
$ws.Range("D4").Value = "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.de-de.xlf"
$ws.Range("E4").Value = "2016-03-21 19:01:30"
$ws.Range("I4").Value = "e2e\calleeMd1.md,`ne2e\calleeMd2.md"
$ws.Range("J4").Value = "Include"
$ws.Range("K4").ClearContents()

$ws.Range("A5").Value = "callerMd2.md"
$ws.Range("B5").Value = ".md"
$ws.Range("C5").Value = "Ready for handoff"
$ws.Range("D5").Value = "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.de-de.xlf"
$ws.Range("E5").Value = "2016-03-21 19:01:30"
$ws.Range("H5").Value = "0001-01-01 00:00:00"
$ws.Range("I5").Value = "e2e\calleeMd1.md"
$ws.Range("J5").Value = "Include"
$ws.Range("E5").NumberFormat = $dateFmt
$ws.Range("H5").NumberFormat = $dateFmt

$ws.Hyperlinks.Delete()
$baseMd = "https://github.com/OpenLocalizationTest/oltest/blob/d21edc54a3494f898cf9a13269ff04c61eedb23a/e2e/"
$baseXlf = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fbcde5b43aa3fc2840cfe4ec712cf025e5220223/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/"
$ws.Hyperlinks.Add($ws.Range("A2"), ($baseMd + "calleeMd1.md"), [Type]::Missing, [Type]::Missing, "calleeMd1.md")
$ws.Hyperlinks.Add($ws.Range("D2"), ($baseXlf + "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.de-de.xlf"), [Type]::Missing, [Type]::Missing, "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), ($baseMd + "calleeMd2.md"), [Type]::Missing, [Type]::Missing, "calleeMd2.md")
$ws.Hyperlinks.Add($ws.Range("D3"), ($baseXlf + "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.de-de.xlf"), [Type]::Missing, [Type]::Missing, "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A4"), ($baseMd + "callerMd1.md"), [Type]::Missing, [Type]::Missing, "callerMd1.md")
$ws.Hyperlinks.Add($ws.Range("D4"), ($baseXlf + "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.de-de.xlf"), [Type]::Missing, [Type]::Missing, "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A5"), ($baseMd + "callerMd2.md"), [Type]::Missing, [Type]::Missing, "callerMd2.md")
$ws.Hyperlinks.Add($ws.Range("D5"), ($baseXlf + "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.de-de.xlf"), [Type]::Missing, [Type]::Missing, "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.de-de.xlf")
